$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ11702297",
    "summ11938865",
    "summ12184397",
    "summ12418467",
    "summ12666321",
    "summ12901386",
    "summ13151380",
    "summ13417839",
    "summ13681772"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}
